$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.08%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.16%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.241"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.01%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07661"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.64%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.9177"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.80%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'2.440"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.97%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1255"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'13.52%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1825"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.09%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09189"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.45%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04283"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.12%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.26%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'1.38%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'-0.10%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D17").Value = "'3.356"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'4.327"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.42%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'7.154"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'9.31%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'1.93%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2898"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'8.08%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04071"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.40%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001264"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.35%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004147"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'3.70%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-2.09%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02461"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'2.74%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05287"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.00%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007856"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.07%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1315"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.06%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006836"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.90%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001912"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.01%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007749"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.16%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3043"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.61%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006724"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.14%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.33%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'1,696.84%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'0.00002108"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.33%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.33%"
$ws.Range("E51").Style = "Normal"
